# Adding teachers to excel file
# Updates the timetable cells on the "Time Table" sheet so each scheduled
# course/slot also carries the assigned teacher initials in brackets
# (e.g. "CS894  /  " becomes "CS894[SMa]  /  []").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Free Period!"
$ws.Range("B4").Value = "Free Period!"
$ws.Range("C4").Value = "Free Period!"
$ws.Range("D4").Value = "Free Period!"
$ws.Range("E4").Value = "CS894[SMa]  /  []"
$ws.Range("F4").Value = "CS894[SMa]  /  []"
$ws.Range("G4").Value = "CS894[SMa]  /  []"

$ws.Range("A6").Value = "Free Period!"
$ws.Range("B6").Value = "CS891[DC, SMa]  /  []"
$ws.Range("C6").Value = "CS891[DC, SMa]  /  []"
$ws.Range("D6").Value = "CS891[DC, SMa]  /  []"
$ws.Range("E6").Value = "CS801A[SSK]  /  []"
$ws.Range("F6").Value = "CS802A[BDu]  /  CS802B[BDu]"
$ws.Range("G6").Value = "HU801[KB]  /  []"

$ws.Range("A8").Value = "Free Period!"
$ws.Range("B8").Value = "CS802A[BDu]  /  CS802B[BDu]"
$ws.Range("C8").Value = "Free Period!"
$ws.Range("D8").Value = "CS801A[SSK]  /  []"
$ws.Range("E8").Value = "CS894[SSK]  /  []"
$ws.Range("F8").Value = "CS894[SSK]  /  []"
$ws.Range("G8").Value = "CS894[SSK]  /  []"

$ws.Range("A10").Value = "CS891[DC, SMa]  /  []"
$ws.Range("B10").Value = "CS891[DC, SMa]  /  []"
$ws.Range("C10").Value = "CS891[DC, SMa]  /  []"
$ws.Range("D10").Value = "CS802A[BDu]  /  CS802B[BDu]"
$ws.Range("E10").Value = "CS894[SMa]  /  []"
$ws.Range("F10").Value = "CS894[SMa]  /  []"
$ws.Range("G10").Value = "CS894[SMa]  /  []"

$ws.Range("A12").Value = "CS894[SMa]  /  []"
$ws.Range("B12").Value = "CS894[SMa]  /  []"
$ws.Range("C12").Value = "CS894[SMa]  /  []"
$ws.Range("D12").Value = "CS801A[SSK]  /  []"
$ws.Range("E12").Value = "Free Period!"
$ws.Range("F12").Value = "Free Period!"
$ws.Range("G12").Value = "HU801[KB]  /  []"

$ws.Range("A16").Value = "IT801A[RCh]  /  []"
$ws.Range("B16").Value = "IT894[AB]  /  []"
$ws.Range("C16").Value = "IT894[AB]  /  []"
$ws.Range("D16").Value = "IT894[AB]  /  []"
$ws.Range("E16").Value = "IT891[KDa, SU]  /  []"
$ws.Range("F16").Value = "IT891[KDa, SU]  /  []"
$ws.Range("G16").Value = "IT891[KDa, SU]  /  []"

$ws.Range("A18").Value = "Free Period!"
$ws.Range("B18").Value = "Free Period!"
$ws.Range("C18").Value = "Free Period!"
$ws.Range("D18").Value = "Free Period!"
$ws.Range("E18").Value = "IT894[AB]  /  []"
$ws.Range("F18").Value = "IT894[AB]  /  []"
$ws.Range("G18").Value = "IT894[AB]  /  []"

$ws.Range("A20").Value = "IT802A[KDa]  /  []"
$ws.Range("B20").Value = "IT894[AKS]  /  []"
$ws.Range("C20").Value = "IT894[AKS]  /  []"
$ws.Range("D20").Value = "IT894[AKS]  /  []"
$ws.Range("E20").Value = "IT891[KDa, SU]  /  []"
$ws.Range("F20").Value = "IT891[KDa, SU]  /  []"
$ws.Range("G20").Value = "IT891[KDa, SU]  /  []"

$ws.Range("A22").Value = "Free Period!"
$ws.Range("B22").Value = "IT801A[RCh]  /  []"
$ws.Range("C22").Value = "Free Period!"
$ws.Range("D22").Value = "Free Period!"
$ws.Range("E22").Value = "IT802A[KDa]  /  []"
$ws.Range("F22").Value = "Free Period!"
$ws.Range("G22").Value = "HU801[KB]  /  []"

$ws.Range("A24").Value = "IT802A[KDa]  /  []"
$ws.Range("B24").Value = "IT801A[RCh]  /  []"
$ws.Range("C24").Value = "HU801[KB]  /  []"
$ws.Range("D24").Value = "Free Period!"
$ws.Range("E24").Value = "IT894[AKS]  /  []"
$ws.Range("F24").Value = "IT894[AKS]  /  []"
$ws.Range("G24").Value = "IT894[AKS]  /  []"

$ws.Range("A28").Value = "HU801[KB]  /  []"
$ws.Range("B28").Value = "ECE894[PC]  /  []"
$ws.Range("C28").Value = "ECE894[PC]  /  []"
$ws.Range("D28").Value = "ECE894[PC]  /  []"
$ws.Range("E28").Value = "ECE891[AnC, BC]  /  []"
$ws.Range("F28").Value = "ECE891[AnC, BC]  /  []"
$ws.Range("G28").Value = "ECE891[AnC, BC]  /  []"

$ws.Range("A30").Value = "ECE802A[ArD]  /  []"
$ws.Range("B30").Value = "Free Period!"
$ws.Range("C30").Value = "Free Period!"
$ws.Range("D30").Value = "Free Period!"
$ws.Range("E30").Value = "ECE801A[DK]  /  []"
$ws.Range("F30").Value = "Free Period!"
$ws.Range("G30").Value = "Free Period!"

$ws.Range("A32").Value = "ECE894[RND]  /  []"
$ws.Range("B32").Value = "ECE894[RND]  /  []"
$ws.Range("C32").Value = "ECE894[RND]  /  []"
$ws.Range("D32").Value = "ECE801A[DK]  /  []"
$ws.Range("E32").Value = "ECE802A[ArD]  /  []"
$ws.Range("F32").Value = "Free Period!"
$ws.Range("G32").Value = "Free Period!"

$ws.Range("A34").Value = "ECE802A[ArD]  /  []"
$ws.Range("B34").Value = "ECE894[PC]  /  []"
$ws.Range("C34").Value = "ECE894[PC]  /  []"
$ws.Range("D34").Value = "ECE894[PC]  /  []"
$ws.Range("E34").Value = "HU801[KB]  /  []"
$ws.Range("F34").Value = "Free Period!"
$ws.Range("G34").Value = "Free Period!"

$ws.Range("A36").Value = "ECE891[SD, BC]  /  []"
$ws.Range("B36").Value = "ECE891[SD, BC]  /  []"
$ws.Range("C36").Value = "ECE891[SD, BC]  /  []"
$ws.Range("D36").Value = "ECE801A[DK]  /  []"
$ws.Range("E36").Value = "ECE894[TD]  /  []"
$ws.Range("F36").Value = "ECE894[TD]  /  []"
$ws.Range("G36").Value = "ECE894[TD]  /  []"

$ws.Range("A40").Value = "EE802A[KR]  /  []"
$ws.Range("B40").Value = "Free Period!"
$ws.Range("C40").Value = "Free Period!"
$ws.Range("D40").Value = "HU801[KB]  /  []"
$ws.Range("E40").Value = "EE894[MB]  /  []"
$ws.Range("F40").Value = "EE894[MB]  /  []"
$ws.Range("G40").Value = "EE894[MB]  /  []"

$ws.Range("A42").Value = "EE801A[ASG]  /  []"
$ws.Range("B42").Value = "EE894[MB]  /  []"
$ws.Range("C42").Value = "EE894[MB]  /  []"
$ws.Range("D42").Value = "EE894[MB]  /  []"
$ws.Range("E42").Value = "EE891[PG, IB]  /  []"
$ws.Range("F42").Value = "EE891[PG, IB]  /  []"
$ws.Range("G42").Value = "EE891[PG, IB]  /  []"

$ws.Range("A44").Value = "EE802A[KR]  /  []"
$ws.Range("B44").Value = "Free Period!"
$ws.Range("C44").Value = "Free Period!"
$ws.Range("D44").Value = "EE801A[ASG]  /  []"
$ws.Range("E44").Value = "EE894[MB]  /  []"
$ws.Range("F44").Value = "EE894[MB]  /  []"
$ws.Range("G44").Value = "EE894[MB]  /  []"

$ws.Range("A46").Value = "EE802A[KR]  /  []"
$ws.Range("B46").Value = "EE801A[ASG]  /  []"
$ws.Range("C46").Value = "Free Period!"
$ws.Range("D46").Value = "HU801[KB]  /  []"
$ws.Range("E46").Value = "Free Period!"
$ws.Range("F46").Value = "Free Period!"
$ws.Range("G46").Value = "Free Period!"

$ws.Range("A48").Value = "Free Period!"
$ws.Range("B48").Value = "EE894[MB]  /  []"
$ws.Range("C48").Value = "EE894[MB]  /  []"
$ws.Range("D48").Value = "EE894[MB]  /  []"
$ws.Range("E48").Value = "EE891[PG, IB]  /  []"
$ws.Range("F48").Value = "EE891[PG, IB]  /  []"
$ws.Range("G48").Value = "EE891[PG, IB]  /  []"
